# Eventbuttons.xlsx - "Commands" sheet update
#
# Adds a new IO Command row documenting the YOCTOPUCE Yocto-Watt
# power-counter reset command:
#     powerReset([sn])  ->  YOCTOPUCE resets the power counter of the Yocto-Watt module
#
# The new row is inserted right after the existing YOCTOPUCE relay "pip(...)"
# row (i.e. becomes row 42), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# Drop the trailing phantom/empty row (present at the very bottom of the
# used range in the source file) before inserting, so the sheet ends up
# with exactly the expected number of rows instead of shifting the
# phantom row further down.
$ws.Rows.Item(1048576).Delete()

# Insert a new blank row at row 42 (shifts old rows 42..129 to 43..130).
$ws.Rows.Item(42).Insert()

# Fill in the new command documentation row.
$ws.Range("B42").Value = "powerReset([sn])"
$ws.Range("C42").Value = "YOCTOPUCE resets the power counter of the Yocto-Watt module"

# Leave the selection where the author ended up after typing the new text.
$ws.Range("C42").Select()
